# Auto-applied numeric updates to the profit-calculation sheets
# (columns H..N: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 2
$ws.Cells.Item(2, 8).Value = 125.90909  # H2: 77.25 -> 125.90909
$ws.Cells.Item(2, 9).Value = 140.66667  # I2: 87.52941 -> 140.66667
$ws.Cells.Item(2, 10).Value = 59.5  # J2: 19 -> 59.5
$ws.Cells.Item(2, 11).Value = 140.66667  # K2: 87.52941 -> 140.66667
$ws.Cells.Item(2, 12).Value = 59.5  # L2: 19 -> 59.5
$ws.Cells.Item(2, 13).Value = -27.66667000000001  # M2: 25.47059 -> -27.66667000000001
$ws.Cells.Item(2, 14).Value = -285.5  # N2: -245 -> -285.5

# ALC!row 4
$ws.Cells.Item(4, 8).Value = 338.46155  # H4: 334.53845 -> 338.46155
$ws.Cells.Item(4, 9).Value = 274.54544  # I4: 232 -> 274.54544
$ws.Cells.Item(4, 10).Value = 690  # J4: 676.3333 -> 690
$ws.Cells.Item(4, 11).Value = 274.54544  # K4: 232 -> 274.54544
$ws.Cells.Item(4, 12).Value = 690  # L4: 676.3333 -> 690
$ws.Cells.Item(4, 13).Value = -160.54544  # M4: -118 -> -160.54544
$ws.Cells.Item(4, 14).Value = -918  # N4: -904.3333 -> -918

# ALC!row 18
$ws.Cells.Item(18, 8).Value = 1480  # H18: 1772.7273 -> 1480
$ws.Cells.Item(18, 9).Value = 1220.1  # I18: 1642.8572 -> 1220.1
$ws.Cells.Item(18, 10).Value = 1999.8  # J18: 2000 -> 1999.8
$ws.Cells.Item(18, 11).Value = 1220.1  # K18: 1642.8572 -> 1220.1
$ws.Cells.Item(18, 12).Value = 1999.8  # L18: 2000 -> 1999.8
$ws.Cells.Item(18, 13).Value = -936.0999999999999  # M18: -1358.8572 -> -936.0999999999999
$ws.Cells.Item(18, 14).Value = -2567.8  # N18: -2568 -> -2567.8

# ALC!row 62
$ws.Cells.Item(62, 8).Value = 27782904  # H62: 55560556 -> 27782904
$ws.Cells.Item(62, 9).Value = 37041036  # I62: 55560556 -> 37041036
$ws.Cells.Item(62, 10).Value = 8500  # J62: 0 -> 8500
$ws.Cells.Item(62, 11).Value = 37041036  # K62: 55560556 -> 37041036
$ws.Cells.Item(62, 12).Value = 8500  # L62: 0 -> 8500
$ws.Cells.Item(62, 13).Value = -37040412  # M62: -55559932 -> -37040412
$ws.Cells.Item(62, 14).Value = -9748  # N62: None -> -9748

# ALC!row 65
$ws.Cells.Item(65, 8).Value = 27782904  # H65: 55560556 -> 27782904
$ws.Cells.Item(65, 9).Value = 37041036  # I65: 55560556 -> 37041036
$ws.Cells.Item(65, 10).Value = 8500  # J65: 0 -> 8500
$ws.Cells.Item(65, 11).Value = 185205180  # K65: 277802780 -> 185205180
$ws.Cells.Item(65, 12).Value = 42500  # L65: 0 -> 42500
$ws.Cells.Item(65, 13).Value = -185202060  # M65: -277799660 -> -185202060
$ws.Cells.Item(65, 14).Value = -48740  # N65: None -> -48740

# ALC!row 70
$ws.Cells.Item(70, 9).Value = 1745.8  # I70: 1621.5 -> 1745.8
$ws.Cells.Item(70, 10).Value = 1598.3  # J70: 1664.7778 -> 1598.3
$ws.Cells.Item(70, 11).Value = 5237.4  # K70: 4864.5 -> 5237.4
$ws.Cells.Item(70, 12).Value = 4794.9  # L70: 4994.3334 -> 4794.9
$ws.Cells.Item(70, 13).Value = -4967.4  # M70: -4594.5 -> -4967.4
$ws.Cells.Item(70, 14).Value = -5334.9  # N70: -5534.3334 -> -5334.9

# ALC!row 73
$ws.Cells.Item(73, 9).Value = 1745.8  # I73: 1621.5 -> 1745.8
$ws.Cells.Item(73, 10).Value = 1598.3  # J73: 1664.7778 -> 1598.3
$ws.Cells.Item(73, 11).Value = 5237.4  # K73: 4864.5 -> 5237.4
$ws.Cells.Item(73, 12).Value = 4794.9  # L73: 4994.3334 -> 4794.9
$ws.Cells.Item(73, 13).Value = -4301.4  # M73: -3928.5 -> -4301.4
$ws.Cells.Item(73, 14).Value = -6666.9  # N73: -6866.3334 -> -6666.9

# ALC!row 112
$ws.Cells.Item(112, 8).Value = 2781.9512  # H112: 2757.561 -> 2781.9512
$ws.Cells.Item(112, 9).Value = 1066.6666  # I112: 1050 -> 1066.6666
$ws.Cells.Item(112, 10).Value = 2917.3684  # J112: 2942.162 -> 2917.3684
$ws.Cells.Item(112, 11).Value = 3199.9998  # K112: 3150 -> 3199.9998
$ws.Cells.Item(112, 12).Value = 8752.1052  # L112: 8826.485999999999 -> 8752.1052
$ws.Cells.Item(112, 13).Value = -2091.9998  # M112: -2042 -> -2091.9998
$ws.Cells.Item(112, 14).Value = -10968.1052  # N112: -11042.486 -> -10968.1052

# ALC!row 137
$ws.Cells.Item(137, 8).Value = 1668.871  # H137: 1668.9032 -> 1668.871
$ws.Cells.Item(137, 9).Value = 1458.5834  # I137: 1458.6666 -> 1458.5834
$ws.Cells.Item(137, 11).Value = 4375.7502  # K137: 4375.9998 -> 4375.7502
$ws.Cells.Item(137, 13).Value = -1825.7502  # M137: -1825.9998 -> -1825.7502

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 74
$ws.Cells.Item(74, 8).Value = 2316.2307  # H74: 2261.5 -> 2316.2307
$ws.Cells.Item(74, 9).Value = 1919.1818  # I74: 2018.6666 -> 1919.1818
$ws.Cells.Item(74, 10).Value = 4500  # J74: 2990 -> 4500
$ws.Cells.Item(74, 11).Value = 1919.1818  # K74: 2018.6666 -> 1919.1818
$ws.Cells.Item(74, 12).Value = 4500  # L74: 2990 -> 4500
$ws.Cells.Item(74, 13).Value = -1045.1818  # M74: -1144.6666 -> -1045.1818
$ws.Cells.Item(74, 14).Value = -6248  # N74: -4738 -> -6248

# ARM!row 77
$ws.Cells.Item(77, 8).Value = 2316.2307  # H77: 2261.5 -> 2316.2307
$ws.Cells.Item(77, 9).Value = 1919.1818  # I77: 2018.6666 -> 1919.1818
$ws.Cells.Item(77, 10).Value = 4500  # J77: 2990 -> 4500
$ws.Cells.Item(77, 11).Value = 9595.909  # K77: 10093.333 -> 9595.909
$ws.Cells.Item(77, 12).Value = 22500  # L77: 14950 -> 22500
$ws.Cells.Item(77, 13).Value = -5227.909  # M77: -5725.333000000001 -> -5227.909
$ws.Cells.Item(77, 14).Value = -31236  # N77: -23686 -> -31236

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 31
$ws.Cells.Item(31, 8).Value = 1047.7843  # H31: 1067.34 -> 1047.7843
$ws.Cells.Item(31, 9).Value = 1003.4146  # I31: 1026.0975 -> 1003.4146
$ws.Cells.Item(31, 10).Value = 1229.7  # J31: 1255.2222 -> 1229.7
$ws.Cells.Item(31, 11).Value = 1003.4146  # K31: 1026.0975 -> 1003.4146
$ws.Cells.Item(31, 12).Value = 1229.7  # L31: 1255.2222 -> 1229.7
$ws.Cells.Item(31, 13).Value = -708.4146  # M31: -731.0975000000001 -> -708.4146
$ws.Cells.Item(31, 14).Value = -1819.7  # N31: -1845.2222 -> -1819.7

# CRP!row 34
$ws.Cells.Item(34, 8).Value = 1047.7843  # H34: 1067.34 -> 1047.7843
$ws.Cells.Item(34, 9).Value = 1003.4146  # I34: 1026.0975 -> 1003.4146
$ws.Cells.Item(34, 10).Value = 1229.7  # J34: 1255.2222 -> 1229.7
$ws.Cells.Item(34, 11).Value = 1003.4146  # K34: 1026.0975 -> 1003.4146
$ws.Cells.Item(34, 12).Value = 1229.7  # L34: 1255.2222 -> 1229.7
$ws.Cells.Item(34, 13).Value = -801.4146  # M34: -824.0975000000001 -> -801.4146
$ws.Cells.Item(34, 14).Value = -1633.7  # N34: -1659.2222 -> -1633.7

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 38
$ws.Cells.Item(38, 8).Value = 40  # H38: 36.666668 -> 40
$ws.Cells.Item(38, 9).Value = 45  # I38: 40 -> 45
$ws.Cells.Item(38, 10).Value = 36.666668  # J38: 35 -> 36.666668
$ws.Cells.Item(38, 11).Value = 135  # K38: 120 -> 135
$ws.Cells.Item(38, 12).Value = 110.000004  # L38: 105 -> 110.000004
$ws.Cells.Item(38, 13).Value = 212  # M38: 227 -> 212
$ws.Cells.Item(38, 14).Value = -804.000004  # N38: -799 -> -804.000004

# CUL!row 40
$ws.Cells.Item(40, 8).Value = 254.18182  # H40: 313.375 -> 254.18182
$ws.Cells.Item(40, 9).Value = 124.5  # I40: 103 -> 124.5
$ws.Cells.Item(40, 10).Value = 328.2857  # J40: 409 -> 328.2857
$ws.Cells.Item(40, 11).Value = 498  # K40: 412 -> 498
$ws.Cells.Item(40, 12).Value = 1313.1428  # L40: 1636 -> 1313.1428
$ws.Cells.Item(40, 13).Value = -429  # M40: -343 -> -429
$ws.Cells.Item(40, 14).Value = -1451.1428  # N40: -1774 -> -1451.1428

# CUL!row 80
$ws.Cells.Item(80, 8).Value = 4332.222  # H80: 4499.5 -> 4332.222
$ws.Cells.Item(80, 9).Value = 0  # I80: 2998 -> 0
$ws.Cells.Item(80, 10).Value = 4332.222  # J80: 5000 -> 4332.222
$ws.Cells.Item(80, 11).Value = 0  # K80: 8994 -> 0
$ws.Cells.Item(80, 12).ClearContents()  # L80: 15000 -> (removed)
$ws.Cells.Item(80, 13).ClearContents()  # M80: -8058 -> (removed)
$ws.Cells.Item(80, 14).Value = -14868.666  # N80: -16872 -> -14868.666

# CUL!row 83
$ws.Cells.Item(83, 8).Value = 4332.222  # H83: 4499.5 -> 4332.222
$ws.Cells.Item(83, 9).Value = 0  # I83: 2998 -> 0
$ws.Cells.Item(83, 10).Value = 4332.222  # J83: 5000 -> 4332.222
$ws.Cells.Item(83, 11).Value = 0  # K83: 26982 -> 0
$ws.Cells.Item(83, 12).ClearContents()  # L83: 45000 -> (removed)
$ws.Cells.Item(83, 13).ClearContents()  # M83: -22302 -> (removed)
$ws.Cells.Item(83, 14).Value = -48349.998  # N83: -54360 -> -48349.998

# CUL!row 87
$ws.Cells.Item(87, 8).Value = 835.6  # H87: 2089.7144 -> 835.6
$ws.Cells.Item(87, 9).Value = 835.6  # I87: 882 -> 835.6
$ws.Cells.Item(87, 10).Value = 0  # J87: 3700 -> 0
$ws.Cells.Item(87, 11).Value = 2506.8  # K87: 2646 -> 2506.8
$ws.Cells.Item(87, 12).Value = 0  # L87: 11100 -> 0
$ws.Cells.Item(87, 13).ClearContents()  # M87: -1398 -> (removed)
$ws.Cells.Item(87, 14).ClearContents()  # N87: -13596 -> (removed)

# CUL!row 90
$ws.Cells.Item(90, 8).Value = 835.6  # H90: 2089.7144 -> 835.6
$ws.Cells.Item(90, 9).Value = 835.6  # I90: 882 -> 835.6
$ws.Cells.Item(90, 10).Value = 0  # J90: 3700 -> 0
$ws.Cells.Item(90, 11).Value = 7520.400000000001  # K90: 7938 -> 7520.400000000001
$ws.Cells.Item(90, 12).Value = 0  # L90: 33300 -> 0
$ws.Cells.Item(90, 13).ClearContents()  # M90: -1698 -> (removed)
$ws.Cells.Item(90, 14).ClearContents()  # N90: -45780 -> (removed)

# CUL!row 131
$ws.Cells.Item(131, 8).Value = 18871190  # H131: 16952166 -> 18871190
$ws.Cells.Item(131, 9).Value = 100000320  # I131: 111111416 -> 100000320
$ws.Cells.Item(131, 10).Value = 3951.3489  # J131: 3500.72 -> 3951.3489
$ws.Cells.Item(131, 11).Value = 300000960  # K131: 333334248 -> 300000960
$ws.Cells.Item(131, 12).Value = 11854.0467  # L131: 10502.16 -> 11854.0467
$ws.Cells.Item(131, 13).Value = -299995920  # M131: -333329208 -> -299995920
$ws.Cells.Item(131, 14).Value = -21934.0467  # N131: -20582.16 -> -21934.0467

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 86
$ws.Cells.Item(86, 8).Value = 23197.6  # H86: 25997 -> 23197.6
$ws.Cells.Item(86, 10).Value = 23197.6  # J86: 25997 -> 23197.6
$ws.Cells.Item(86, 12).Value = 23197.6  # L86: 25997 -> 23197.6
$ws.Cells.Item(86, 14).Value = -25569.6  # N86: -28369 -> -25569.6

# GSM!row 89
$ws.Cells.Item(89, 8).Value = 23197.6  # H89: 25997 -> 23197.6
$ws.Cells.Item(89, 10).Value = 23197.6  # J89: 25997 -> 23197.6
$ws.Cells.Item(89, 12).Value = 69592.79999999999  # L89: 77991 -> 69592.79999999999
$ws.Cells.Item(89, 14).Value = -81448.79999999999  # N89: -89847 -> -81448.79999999999

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 22
$ws.Cells.Item(22, 8).Value = 1041.5714  # H22: 966.44446 -> 1041.5714
$ws.Cells.Item(22, 10).Value = 1500  # J22: 1101.75 -> 1500
$ws.Cells.Item(22, 12).Value = 1500  # L22: 1101.75 -> 1500
$ws.Cells.Item(22, 14).Value = -2090  # N22: -1691.75 -> -2090

# LTW!row 27
$ws.Cells.Item(27, 8).Value = 1041.5714  # H27: 966.44446 -> 1041.5714
$ws.Cells.Item(27, 10).Value = 1500  # J27: 1101.75 -> 1500
$ws.Cells.Item(27, 12).Value = 1500  # L27: 1101.75 -> 1500
$ws.Cells.Item(27, 14).Value = -1714  # N27: -1315.75 -> -1714

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 81
$ws.Cells.Item(81, 8).Value = 0  # H81: 2000 -> 0
$ws.Cells.Item(81, 10).Value = 0  # J81: 2000 -> 0
$ws.Cells.Item(81, 12).ClearContents()  # L81: 4000 -> (removed)
$ws.Cells.Item(81, 14).ClearContents()  # N81: -6122 -> (removed)

# WVR!row 84
$ws.Cells.Item(84, 8).Value = 0  # H84: 2000 -> 0
$ws.Cells.Item(84, 10).Value = 0  # J84: 2000 -> 0
$ws.Cells.Item(84, 12).ClearContents()  # L84: 20000 -> (removed)
$ws.Cells.Item(84, 14).ClearContents()  # N84: -30608 -> (removed)
